$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: values are prefixed with a leading apostrophe to force
# Excel to store them as literal text (matching the original inline-string
# cell type) rather than auto-converting numeric-looking text to numbers.
$ws.Range("D2").Value = "'67.758.46"
$ws.Range("E2").Value = "'  +0.13%  "
$ws.Range("D3").Value = "'3.804.53"
$ws.Range("E3").Value = "'  +0.26%  "
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("D5").Value = "'604.60"
$ws.Range("E5").Value = "'  +1.59%  "
$ws.Range("D6").Value = "'167.40"
$ws.Range("E6").Value = "'  +0.47%  "
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E9").Value = "'  +0.83%  "
$ws.Range("E10").Value = "'  -0.68%  "
$ws.Range("E11").Value = "'  -0.03%  "
$ws.Range("E12").Value = "'  -1.09%  "
$ws.Range("D13").Value = "'35.98"
$ws.Range("E13").Value = "'  -1.04%  "
$ws.Range("D14").Value = "'4.442.29"
$ws.Range("E14").Value = "'  +0.22%  "
$ws.Range("D15").Value = "'3.783.29"
$ws.Range("E15").Value = "'  -0.35%  "
$ws.Range("D16").Value = "'18.50"
$ws.Range("E16").Value = "'  -0.80%  "
$ws.Range("D17").Value = "'67.844.68"
$ws.Range("E17").Value = "'  +0.25%  "
$ws.Range("D18").Value = "'7.08"
$ws.Range("E18").Value = "'  +1.30%  "
$ws.Range("E19").Value = "'  +0.50%  "
$ws.Range("D20").Value = "'462.29"
$ws.Range("E20").Value = "'  +0.93%  "
$ws.Range("E21").Value = "'  -3.53%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "'  +0.43%  "
$ws.Range("E23").Value = "'  -2.28%  "
$ws.Range("D24").Value = "'83.36"
$ws.Range("E24").Value = "'  -0.24%  "
$ws.Range("D25").Value = "'12.10"
$ws.Range("D26").Value = "'2.10"
$ws.Range("E26").Value = "'  -1.36%  "
$ws.Range("B27").Value = "'Dai"
$ws.Range("C27").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "'  +0.10%  "
$ws.Range("B28").Value = "'RenderToken"
$ws.Range("C28").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.02"
$ws.Range("E28").Value = "'  -0.66%  "
$ws.Range("D29").Value = "'3.953.33"
$ws.Range("E29").Value = "'  +0.23%  "
$ws.Range("D30").Value = "'2.79"
$ws.Range("E30").Value = "'  +0.03%  "
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "'  +1.69%  "
$ws.Range("E32").Value = "'  +1.64%  "
$ws.Range("D33").Value = "'29.57"
$ws.Range("E33").Value = "'  -0.84%  "
$ws.Range("E34").Value = "'  +0.04%  "
$ws.Range("E35").Value = "'  -1.29%  "
$ws.Range("D36").Value = "'3.744.93"
$ws.Range("E36").Value = "'  -0.10%  "
$ws.Range("E37").Value = "'  -0.16%  "
$ws.Range("D38").Value = "'3.40"
$ws.Range("E38").Value = "'  +1.60%  "
$ws.Range("E39").Value = "'  -0.06%  "
$ws.Range("E40").Value = "'  +0.10%  "
$ws.Range("E41").Value = "'  +0.44%  "
$ws.Range("E42").Value = "'  -0.07%  "
$ws.Range("E43").Value = "'  -0.02%  "
$ws.Range("E44").Value = "'  +2.14%  "
$ws.Range("E45").Value = "'  +0.71%  "
$ws.Range("D46").Value = "'42.96"
$ws.Range("E46").Value = "'  -3.87%  "
$ws.Range("D47").Value = "'27.79"
$ws.Range("E47").Value = "'  +9.91%  "
$ws.Range("E48").Value = "'  -0.42%  "
$ws.Range("E49").Value = "'  +9.50%  "
$ws.Range("D50").Value = "'148.40"
$ws.Range("E50").Value = "'  -0.40%  "
$ws.Range("D51").Value = "'1.84"
$ws.Range("E51").Value = "'  +0.05%  "
